# The deck currently has the "Integral" theme (Red Violet colour scheme)
# applied to its slide master. This script switches the applied theme's
# colour scheme back to the standard default "Office Theme" colours
# (the palette that was otherwise sitting unused in the package), which
# is the effective, user-visible result of the underlying edit.
#
# PowerPoint exposes the twelve theme colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) through Master.Theme.ThemeColorScheme, in
# that fixed order - exactly mirroring <a:clrScheme> in the theme part.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# RGB() isn't available in this host, so pack R+G*256+B*65536 ourselves
# (standard OLE COLORREF / VBA RGB() encoding).
function Pack-Rgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock Office theme colours.
$colorScheme.Item(1).RGB  = Pack-Rgb 0x00 0x00 0x00   # dk1
$colorScheme.Item(2).RGB  = Pack-Rgb 0xFF 0xFF 0xFF   # lt1
$colorScheme.Item(3).RGB  = Pack-Rgb 0x44 0x54 0x6A   # dk2
$colorScheme.Item(4).RGB  = Pack-Rgb 0xE7 0xE6 0xE6   # lt2
$colorScheme.Item(5).RGB  = Pack-Rgb 0x5B 0x9B 0xD5   # accent1
$colorScheme.Item(6).RGB  = Pack-Rgb 0xED 0x7D 0x31   # accent2
$colorScheme.Item(7).RGB  = Pack-Rgb 0xA5 0xA5 0xA5   # accent3
$colorScheme.Item(8).RGB  = Pack-Rgb 0xFF 0xC0 0x00   # accent4
$colorScheme.Item(9).RGB  = Pack-Rgb 0x44 0x72 0xC4   # accent5
$colorScheme.Item(10).RGB = Pack-Rgb 0x70 0xAD 0x47   # accent6
$colorScheme.Item(11).RGB = Pack-Rgb 0x05 0x63 0xC1   # hlink
$colorScheme.Item(12).RGB = Pack-Rgb 0x95 0x4F 0x72   # folHlink
